$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(36)
$p2 = $d.Paragraphs.Item(53)
$target = $d.Range($p1.Range.Start, $p2.Range.End)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="002A084D" w:rsidRPr="00FF5162" w:rsidRDefault="005C06E3"><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="sv-SE"/></w:rPr></w:pPr><w:r w:rsidRPr="00FF5162"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t xml:space="preserve">Rasmus </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00FF5162"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t>Tilljander</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00FF5162"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t xml:space="preserve"> - rati10@student.bth.se</w:t></w:r><w:r w:rsidR="002A084D" w:rsidRPr="00FF5162"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:br/></w:r><w:r w:rsidRPr="00FF5162"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t>Nils Forsman - nifo08@student.bth.se</w:t></w:r><w:r w:rsidR="002A084D" w:rsidRPr="00FF5162"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:br/></w:r><w:r w:rsidRPr="00FF5162"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t>Calle Ketola - cake10@student.bth.se</w:t></w:r><w:r w:rsidR="002A084D" w:rsidRPr="00FF5162"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:br/></w:r><w:r w:rsidRPr="00FF5162"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t>Kim Hansson - kiha10@student.bth.</w:t></w:r><w:r w:rsidR="00E75D8C" w:rsidRPr="00FF5162"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t>se</w:t></w:r><w:r w:rsidR="002A084D" w:rsidRPr="00FF5162"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Introduction</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">In this assignment we were to evaluate a premade architecture of a forest harvesting machine control system called Blunderjack using the formal architecture evaluation AADL with the program OSATE. We were to complete the AADL model after certain specifications and then evaluate it using two specific scenarios we were given. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Evaluation result</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>. the changes made to the model</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="1304"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>Each subsystem now has its own budget for MIPS and RAM resources, based upon the specification given in the assignment.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="1304"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> The predefined flow paths for Device</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>Manager and Safety</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>Manager have been connected successfully.</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> Furthermore,</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">we added a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>rotation_inout_flow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> to the Safety Manager for usage in scenario 2.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="1304"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>In the system configuration we created 3 different flow paths, the first one (</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>F</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>1) for scenario 1 and the last two (</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>F</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">2 &amp; </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>F</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>3) for scenario 2.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="1304"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>We assigned the subsystems to processor and their memory banks</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> based upon configuration 2 &amp; 3</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> as can be seen in our AADL model</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>. Based upon the specification we concluded that we co</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>uld not assign the subsystems in configuration 1.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>2. the result of the evaluations</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="1304"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>Our result for scenario 1</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> the highest late</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">ncy </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>was</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>265</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">ms, this was well below the requirement given. Our conclusion is that given the current system architecture the system fulfils the requirements. </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="1304"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>The result for scenario 2 had a latency</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> of 435ms</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> with</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> the</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> asynchronous </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>test and a latency of 345ms with the synchronous test. Our conclusion is that for the system to fulfil the requirements, the system has to be run synchronous.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)
Write-Output ("ParaCount=" + $d.Paragraphs.Count)
